$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing data-row style (from D8, which already carries it) onto
# all the new cells before filling them in, so they match the rest of the
# table's look (Consolas 8pt green, vertical-centered).
$ws.Range("D8").Copy()
$ws.Range("E8:F9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D10").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("A16:B18").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I2:J2").PasteSpecial(-4122)   # xlPasteFormats (stray formatting near new pasted area)

# New parotid / spinal cord naming variants appended to existing rows 8-10
$ws.Range("E8").Value = "parotid_l"
$ws.Range("F8").Value = "parotid_r"
$ws.Range("D10").Value = "spinalcord"

# New rows of CTV naming variants (GSTT data)
$ws.Range("A16").Value = "CTV60"
$ws.Range("B16").Value = "CTV54combi"

$ws.Range("E9").Value = "Parotid_L"
$ws.Range("F9").Value = "Parotid_R"

$ws.Range("A17").Value = "CTVp"
$ws.Range("B17").Value = "CTVn"
$ws.Range("A18").Value = "CTV1"
$ws.Range("B18").Value = "CTV2"

# Column I sized to fit its (now cleared/unused) pasted content
$ws.Columns("I").ColumnWidth = 11.29

# Leftover UI selection state from editing session
$ws.Range("I1:K5").Select()
